# Quarterly indexing esoteric bug-fix operation
# For every date in column A (rows 2-73), shift it forward by one month
# and pin the day-of-month to the 15th (mid-month quarter index date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date epoch (serial 0 == 1899-12-30)
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$lastRow = 73
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $serial = $cell.Value2

    if ($serial -eq $null) {
        continue
    }

    $d = $epoch.AddDays($serial)
    $d = $d.AddMonths(1)

    $newDate = Get-Date -Year $d.Year -Month $d.Month -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

    $cell.Value2 = $newDate
}
